# Generate Report for Handback
# ------------------------------------------------------------------
# The handback run completed successfully: the "Status" now reads
# "Handed back: in sync with en-US" instead of "Ready for handoff",
# the Latest Handback DateTime stamps advance to the new handback
# run, the stale "handback file is not the latest" Error Detail is
# cleared out, and the affected Status/Latest-Handback-DateTime
# columns are widened to fit the new text.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------- Overview sheet ----------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Widen columns E and F so the longer status text fits.
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------- zh-cn sheet ----------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-24 18:58:44"
$zhcn.Range("P2").Value = ""

# Widen the Status column, shrink the now-empty Error Detail column.
$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(16).ColumnWidth = 12.85

# ---------------- de-de sheet ----------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-24 18:58:51"
$dede.Range("P2").Value = ""

# Widen the Status column, shrink the now-empty Error Detail column.
$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(16).ColumnWidth = 12.85
